# Update gh-pages to output generated at 456a3b4
#
# Both the "展览" (sheet 1) and "全部类型" (sheet 4) tabs gained a new row:
#   合肥·WA二次元饭局, inserted right before the existing
#   "合肥·梦时空SPO1动漫展（取消）" row (i.e. at worksheet row 13), pushing every
#   row from there on down by one. A handful of the "想去人数" (F column)
#   counters were also refreshed (site scrape counts moved up a little).

$wb = $excel.ActiveWorkbook

function Update-FanExpoSheet {
    param($ws)

    # --- refresh "want to go" counters (F column) for the rows above the insert point ---
    $ws.Range("F2").Value  = 241
    $ws.Range("F3").Value  = 267
    $ws.Range("F4").Value  = 278
    $ws.Range("F5").Value  = 822
    $ws.Range("F6").Value  = 272
    $ws.Range("F7").Value  = 6606
    $ws.Range("F8").Value  = 54
    $ws.Range("F9").Value  = 73
    $ws.Range("F10").Value = 115
    $ws.Range("F12").Value = 36

    # --- insert the new event row at row 13, shifting everything else down ---
    $ws.Range("A13").EntireRow.Insert()

    # copy row 14's formatting (bold/border/center-top) down onto the fresh row 13
    $ws.Range("A14").Copy()
    $ws.Range("A13").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    $ws.Range("A13").Value = 12
    $ws.Range("B13").Value = "2024-05-18"
    $ws.Range("C13").Value = "合肥·WA二次元饭局"
    $ws.Range("D13").Value = "临泉路胜利路交叉路（中环国际大厦对面） 太太满庭芳(胜利路店)"
    $ws.Range("E13").Value = "2024.05.18 14:50-05.18 20:00"
    $ws.Range("F13").Value = 8
    $ws.Range("G13").Value = 118
    $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83978"
    $ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202404/wK9Yq9Ta1712657384067.jpeg"

    # --- refresh counters on the rows that shifted down (their own scrape counts moved too) ---
    $ws.Range("F15").Value = 16
    $ws.Range("F16").Value = 212
    $ws.Range("F17").Value = 544
}

$ws1 = $wb.Worksheets.Item(1)   # 展览
Update-FanExpoSheet $ws1

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
Update-FanExpoSheet $ws4
